$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("InputsOptional")
$ws.Columns.Item(2).Insert()
$ws.Cells.Item(1,2).Value = "clip"
$src = $ws.Cells.Item(3, 1)
$src.Copy()
$ws.Cells.Item(3, 2).PasteSpecial(-4163)
$ws.Columns.Item(2).AutoFit()
